$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws_ALC.Range("H17").Value = 3843.8
$ws_ALC.Range("J17").Value = 3940.842
$ws_ALC.Range("L17").Value = 11822.526
$ws_ALC.Range("N17").Value = -12158.526

# ALC row 51
$ws_ALC.Range("H51").Value = 11767.4
$ws_ALC.Range("I51").Value = 3000
$ws_ALC.Range("K51").Value = 3000
$ws_ALC.Range("M51").Value = -2516

# ALC row 58
$ws_ALC.Range("H58").Value = 1034.1666
$ws_ALC.Range("J58").Value = 1989
$ws_ALC.Range("L58").Value = 5967
$ws_ALC.Range("N58").Value = -6267

# ALC row 70
$ws_ALC.Range("H70").Value = 5698.933
$ws_ALC.Range("I70").Value = 3648.3333
$ws_ALC.Range("J70").Value = 7066
$ws_ALC.Range("K70").Value = 10944.9999
$ws_ALC.Range("L70").Value = 21198
$ws_ALC.Range("M70").Value = -10674.9999
$ws_ALC.Range("N70").Value = -21738

# ALC row 73
$ws_ALC.Range("H73").Value = 5698.933
$ws_ALC.Range("I73").Value = 3648.3333
$ws_ALC.Range("J73").Value = 7066
$ws_ALC.Range("K73").Value = 10944.9999
$ws_ALC.Range("L73").Value = 21198
$ws_ALC.Range("M73").Value = -10008.9999
$ws_ALC.Range("N73").Value = -23070

# ALC row 101
$ws_ALC.Range("H101").Value = 16668120
$ws_ALC.Range("I101").Value = 33333998
$ws_ALC.Range("K101").Value = 100001994
$ws_ALC.Range("M101").Value = -100000372

# ALC row 113
$ws_ALC.Range("H113").Value = 4151.6665
$ws_ALC.Range("I113").Value = 3977.5
$ws_ALC.Range("K113").Value = 3977.5
$ws_ALC.Range("M113").Value = -723.5

# ALC row 131
$ws_ALC.Range("H131").Value = 9925.166999999999
$ws_ALC.Range("J131").Value = 17750
$ws_ALC.Range("L131").Value = 53250
$ws_ALC.Range("N131").Value = -63330

# ALC row 135
$ws_ALC.Range("H135").Value = 647.5
$ws_ALC.Range("I135").Value = 454.2857
$ws_ALC.Range("K135").Value = 4088.5713
$ws_ALC.Range("M135").Value = -1553.5713

# ALC row 137
$ws_ALC.Range("H137").Value = 1323.25
$ws_ALC.Range("I137").Value = 1247.9
$ws_ALC.Range("J137").Value = 1700
$ws_ALC.Range("K137").Value = 3743.7
$ws_ALC.Range("L137").Value = 5100
$ws_ALC.Range("M137").Value = -1193.7
$ws_ALC.Range("N137").Value = -10200

# ARM row 32
$ws_ARM.Range("H32").Value = 9258.888999999999
$ws_ARM.Range("I32").Value = 9791.25
$ws_ARM.Range("K32").Value = 9791.25
$ws_ARM.Range("M32").Value = -9504.25

# ARM row 61
$ws_ARM.Range("H61").Value = 3486.6296
$ws_ARM.Range("I61").Value = 1338.9
$ws_ARM.Range("J61").Value = 4750
$ws_ARM.Range("K61").Value = 1338.9
$ws_ARM.Range("L61").Value = 4750
$ws_ARM.Range("M61").Value = -1126.9
$ws_ARM.Range("N61").Value = -5174

# ARM row 122
$ws_ARM.Range("H122").Value = 1156.2858
$ws_ARM.Range("I122").Value = 974
$ws_ARM.Range("K122").Value = 2922
$ws_ARM.Range("M122").Value = -472

# ARM row 136
$ws_ARM.Range("H136").Value = 3486.6296
$ws_ARM.Range("I136").Value = 1338.9
$ws_ARM.Range("J136").Value = 4750
$ws_ARM.Range("K136").Value = 4016.7
$ws_ARM.Range("L136").Value = 14250
$ws_ARM.Range("M136").Value = -1466.7
$ws_ARM.Range("N136").Value = -19350

# BSM row 20
$ws_BSM.Range("H20").Value = 286.5
$ws_BSM.Range("I20").Value = 286.5
$ws_BSM.Range("K20").Value = 286.5
$ws_BSM.Range("M20").Value = -39.5

# CRP row 31
$ws_CRP.Range("H31").Value = 2048.25
$ws_CRP.Range("I31").Value = 1562.5
$ws_CRP.Range("J31").Value = 2291.125
$ws_CRP.Range("K31").Value = 1562.5
$ws_CRP.Range("L31").Value = 2291.125
$ws_CRP.Range("M31").Value = -1267.5
$ws_CRP.Range("N31").Value = -2881.125

# CRP row 34
$ws_CRP.Range("H34").Value = 2048.25
$ws_CRP.Range("I34").Value = 1562.5
$ws_CRP.Range("J34").Value = 2291.125
$ws_CRP.Range("K34").Value = 1562.5
$ws_CRP.Range("L34").Value = 2291.125
$ws_CRP.Range("M34").Value = -1360.5
$ws_CRP.Range("N34").Value = -2695.125

# CRP row 58
$ws_CRP.Range("H58").Value = 1795.6
$ws_CRP.Range("I58").Value = 1606.2778
$ws_CRP.Range("J58").Value = 3499.5
$ws_CRP.Range("K58").Value = 1606.2778
$ws_CRP.Range("L58").Value = 3499.5
$ws_CRP.Range("M58").Value = -1403.2778
$ws_CRP.Range("N58").Value = -3905.5

# CRP row 99
$ws_CRP.Range("H99").Value = 3388.0625
$ws_CRP.Range("I99").Value = 2535.625
$ws_CRP.Range("K99").Value = 2535.625
$ws_CRP.Range("M99").Value = -1037.625

# CRP row 126
$ws_CRP.Range("H126").Value = 3388.0625
$ws_CRP.Range("I126").Value = 2535.625
$ws_CRP.Range("K126").Value = 7606.875
$ws_CRP.Range("M126").Value = -5136.875

# CRP row 136
$ws_CRP.Range("H136").Value = 1795.6
$ws_CRP.Range("I136").Value = 1606.2778
$ws_CRP.Range("J136").Value = 3499.5
$ws_CRP.Range("K136").Value = 4818.8334
$ws_CRP.Range("L136").Value = 10498.5
$ws_CRP.Range("M136").Value = -2268.8334
$ws_CRP.Range("N136").Value = -15598.5

# CUL row 4
$ws_CUL.Range("H4").Value = 8461608
$ws_CUL.Range("I4").Value = 8461608
$ws_CUL.Range("K4").Value = 25384824
$ws_CUL.Range("M4").Value = -25384712

# CUL row 11
$ws_CUL.Range("H11").Value = 26804906
$ws_CUL.Range("I11").Value = 31272274
$ws_CUL.Range("K11").Value = 93816822
$ws_CUL.Range("M11").Value = -93816682

# CUL row 34
$ws_CUL.Range("H34").Value = 2233.6667
$ws_CUL.Range("I34").Value = 300
$ws_CUL.Range("J34").Value = 3200.5
$ws_CUL.Range("K34").Value = 900
$ws_CUL.Range("L34").Value = 9601.5
$ws_CUL.Range("M34").Value = -816
$ws_CUL.Range("N34").Value = -9769.5

# CUL row 42
$ws_CUL.Range("H42").Value = 6499.5
$ws_CUL.Range("I42").Value = 6499
$ws_CUL.Range("K42").Value = 19497
$ws_CUL.Range("M42").Value = -18963

# CUL row 103
$ws_CUL.Range("H103").Value = 2959.3333

# CUL row 113
$ws_CUL.Range("H113").Value = 1518.6666
$ws_CUL.Range("I113").Value = 675
$ws_CUL.Range("J113").Value = 1940.5
$ws_CUL.Range("K113").Value = 2025
$ws_CUL.Range("L113").Value = 5821.5
$ws_CUL.Range("M113").Value = 145
$ws_CUL.Range("N113").Value = -10161.5

# CUL row 129
$ws_CUL.Range("H129").Value = 835823.2
$ws_CUL.Range("I129").Value = 1300
$ws_CUL.Range("J129").Value = 1253084.8
$ws_CUL.Range("K129").Value = 3900
$ws_CUL.Range("L129").Value = 3759254.4
$ws_CUL.Range("M129").Value = 1100
$ws_CUL.Range("N129").Value = -3769254.4

# CUL row 132
$ws_CUL.Range("H132").Value = 2763.8333
$ws_CUL.Range("I132").Value = 1289.5
$ws_CUL.Range("K132").Value = 11605.5
$ws_CUL.Range("M132").Value = -9075.5

# GSM row 2
$ws_GSM.Range("H2").Value = 121.28571
$ws_GSM.Range("I2").Value = 149
$ws_GSM.Range("J2").Value = 19.666666
$ws_GSM.Range("K2").Value = 149
$ws_GSM.Range("L2").Value = 19.666666
$ws_GSM.Range("M2").Value = -36
$ws_GSM.Range("N2").Value = -245.666666

# GSM row 11
$ws_GSM.Range("H11").Value = 4458170
$ws_GSM.Range("I11").Value = 5938893.5
$ws_GSM.Range("J11").Value = 15999
$ws_GSM.Range("K11").Value = 5938893.5
$ws_GSM.Range("L11").Value = 15999
$ws_GSM.Range("M11").Value = -5938754.5
$ws_GSM.Range("N11").Value = -16277

# GSM row 31
$ws_GSM.Range("H31").Value = 1250
$ws_GSM.Range("I31").Value = 1250
$ws_GSM.Range("K31").Value = 1250
$ws_GSM.Range("M31").Value = -958

# GSM row 37
$ws_GSM.Range("H37").Value = 1250
$ws_GSM.Range("I37").Value = 1250
$ws_GSM.Range("K37").Value = 1250
$ws_GSM.Range("M37").Value = -973

# GSM row 46
$ws_GSM.Range("H46").Value = 2636.9
$ws_GSM.Range("I46").Value = 2636.9
$ws_GSM.Range("K46").Value = 2636.9
$ws_GSM.Range("M46").Value = -2480.9

# GSM row 57
$ws_GSM.Range("H57").Value = 27791
$ws_GSM.Range("I57").Value = 27791
$ws_GSM.Range("K57").Value = 27791
$ws_GSM.Range("M57").Value = -26971

# GSM row 102
$ws_GSM.Range("H102").Value = 958.5333000000001
$ws_GSM.Range("I102").Value = 812.7143
$ws_GSM.Range("K102").Value = 812.7143
$ws_GSM.Range("M102").Value = 809.2857

# GSM row 122
$ws_GSM.Range("H122").Value = 5900.5
$ws_GSM.Range("I122").Value = 3851
$ws_GSM.Range("K122").Value = 11553
$ws_GSM.Range("M122").Value = -9103

# GSM row 134
$ws_GSM.Range("H134").Value = 500000
$ws_GSM.Range("J134").Value = 500000
$ws_GSM.Range("L134").Value = 1500000
$ws_GSM.Range("N134").Value = -1505070

# LTW row 63
$ws_LTW.Range("H63").Value = 0
$ws_LTW.Range("I63").Value = 0
$ws_LTW.Range("J63").Value = 0
$ws_LTW.Range("K63").Value = 0
$ws_LTW.Range("L63").ClearContents()
$ws_LTW.Range("M63").ClearContents()
$ws_LTW.Range("N63").Value = 0

# LTW row 66
$ws_LTW.Range("H66").Value = 0
$ws_LTW.Range("I66").Value = 0
$ws_LTW.Range("J66").Value = 0
$ws_LTW.Range("K66").Value = 0
$ws_LTW.Range("L66").ClearContents()
$ws_LTW.Range("M66").ClearContents()
$ws_LTW.Range("N66").Value = 0

# LTW row 74
$ws_LTW.Range("H74").Value = 90000
$ws_LTW.Range("I74").Value = 90000
$ws_LTW.Range("K74").Value = 90000
$ws_LTW.Range("M74").Value = -89002

# LTW row 77
$ws_LTW.Range("H77").Value = 90000
$ws_LTW.Range("I77").Value = 90000
$ws_LTW.Range("K77").Value = 270000
$ws_LTW.Range("M77").Value = -265008

# LTW row 122
$ws_LTW.Range("H122").Value = 5784.2
$ws_LTW.Range("J122").Value = 5390.091
$ws_LTW.Range("L122").Value = 16170.273
$ws_LTW.Range("N122").Value = -21070.273

# WVR row 69
$ws_WVR.Range("H69").Value = 0
$ws_WVR.Range("J69").Value = 0
$ws_WVR.Range("L69").ClearContents()
$ws_WVR.Range("N69").Value = 0

# WVR row 72
$ws_WVR.Range("H72").Value = 0
$ws_WVR.Range("J72").Value = 0
$ws_WVR.Range("L72").ClearContents()
$ws_WVR.Range("N72").Value = 0

# WVR row 81
$ws_WVR.Range("H81").Value = 1002449.1
$ws_WVR.Range("I81").Value = 1059.4286
$ws_WVR.Range("J81").Value = 3339025
$ws_WVR.Range("K81").Value = 2118.8572
$ws_WVR.Range("L81").Value = 6678050
$ws_WVR.Range("M81").Value = -1057.8572
$ws_WVR.Range("N81").Value = -6680172

# WVR row 84
$ws_WVR.Range("H84").Value = 1002449.1
$ws_WVR.Range("I84").Value = 1059.4286
$ws_WVR.Range("J84").Value = 3339025
$ws_WVR.Range("K84").Value = 10594.286
$ws_WVR.Range("L84").Value = 33390250
$ws_WVR.Range("M84").Value = -5290.286
$ws_WVR.Range("N84").Value = -33400858

# WVR row 122
$ws_WVR.Range("H122").Value = 616
$ws_WVR.Range("I122").Value = 661.5
$ws_WVR.Range("J122").Value = 525
$ws_WVR.Range("K122").Value = 1984.5
$ws_WVR.Range("L122").Value = 1575
$ws_WVR.Range("M122").Value = 465.5
$ws_WVR.Range("N122").Value = -6475

# WVR row 126
$ws_WVR.Range("H126").Value = 3993.158
$ws_WVR.Range("I126").Value = 2222.1667
$ws_WVR.Range("K126").Value = 6666.500100000001
$ws_WVR.Range("M126").Value = -4196.500100000001
